# Apply the data-cleaning fixes described in the commit:
# - rename header row to snake_case English field names
# - capitalize lowercase connector words ("de", "del", "el") in place names
# - drop the trailing footer/metadata rows (125-129)
# - shrink the used range dimension to A1:D123

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Normalize lowercase "de"/"del"/"el" to capitalized forms in a few cells ---
$ws.Range("B14").Value = "Hidalgo Del Parral"
$ws.Range("A24").Value = "Ciudad De México"
$ws.Range("A37").Value = "Estado De México"
$ws.Range("B45").Value = "Cuetzala Del Progreso"
$ws.Range("B48").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B53").Value = "Encarnación De Díaz"
$ws.Range("B58").Value = "Zapotlán El Grande"
$ws.Range("B81").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B84").Value = "Amealco De Bonfil"
$ws.Range("B86").Value = "San Juan Del Río"
$ws.Range("B89").Value = "Santa María Del Río"
$ws.Range("B109").Value = "Poza Rica De Hidalgo"

# --- Remove the trailing footer/metadata rows (125-129) ---
$ws.Range("A125:D129").EntireRow.Delete()
